# This workbook's data rows (2-17) represent weekly price records that need
# to be re-sorted/reshuffled. Each record (columns A-T) stays internally
# consistent; only the row position changes. Below is the mapping from the
# destination row number to the source row number (both in the original
# sheet, rows 2..17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$destToSrc = @{
    2  = 17
    3  = 14
    4  = 15
    5  = 4
    6  = 7
    7  = 9
    8  = 13
    9  = 12
    10 = 2
    11 = 3
    12 = 10
    13 = 5
    14 = 16
    15 = 11
    16 = 8
    17 = 6
}

$srcRange = $ws.Range("A2:T17")
$srcValues = $srcRange.Value()

$nRows = 16
$nCols = 20

$newValues = New-Object 'object[,]' $nRows, $nCols

for ($destRow = 2; $destRow -le 17; $destRow++) {
    $srcRow = $destToSrc[$destRow]
    $srcIdx = $srcRow - 1   # 1-based row index inside $srcValues (rows start at A2 => index 1)
    $destIdx = $destRow - 2 # 0-based row index inside $newValues

    for ($col = 1; $col -le $nCols; $col++) {
        $newValues[$destIdx, ($col - 1)] = $srcValues[$srcIdx, $col]
    }
}

$destRange = $ws.Range("A2:T17")
$destRange.Value = $newValues
